# Update "想去人数" (F column) counts on the "展览" (sheet 1) and
# "全部类型" (sheet 4) worksheets, matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# Row -> new F value, for the "展览" worksheet (sheet index 1)
$sheet1Updates = @{
    2  = 291
    3  = 1177
    4  = 16701
    5  = 26
    7  = 61
    8  = 3
    9  = 368
    10 = 211
    12 = 11601
    14 = 1275
    15 = 4589
    16 = 421
    18 = 63
    19 = 883
}

# Row -> new F value, for the "全部类型" worksheet (sheet index 4)
$sheet4Updates = @{
    2  = 291
    4  = 1177
    5  = 16701
    6  = 26
    8  = 61
    9  = 3
    10 = 368
    11 = 211
    15 = 11601
    17 = 1275
    18 = 4589
    19 = 421
    21 = 63
    22 = 883
}

$ws1 = $wb.Worksheets.Item(1)
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item(4)
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
